$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value looks like a plain number but must remain text ---
# Force text format first so Excel does not auto-convert them to numeric cells,
# then restore the default "Normal" style so no stray formatting is introduced.
$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D12",
    "D14",
    "D18",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D27",
    "D32",
    "D36",
    "D37",
    "D38",
    "D41",
    "D43",
    "D47",
    "D50",
)
foreach ($cell in $textForceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D5").Value = '579.82'
$ws.Range("D6").Value = '174.79'
$ws.Range("D7").Value = '0.999'
$ws.Range("D12").Value = '0.483'
$ws.Range("D14").Value = '37.52'
$ws.Range("D18").Value = '7.18'
$ws.Range("D20").Value = '16.15'
$ws.Range("D21").Value = '488.42'
$ws.Range("D22").Value = '0.716'
$ws.Range("D24").Value = '84.25'
$ws.Range("D25").Value = '13.26'
$ws.Range("D27").Value = '10.06'
$ws.Range("D32").Value = '28.87'
$ws.Range("D36").Value = '5.93'
$ws.Range("D37").Value = '0.990'
$ws.Range("D38").Value = '47.81'
$ws.Range("D41").Value = '0.312'
$ws.Range("D43").Value = '8.66'
$ws.Range("D47").Value = '383.54'
$ws.Range("D50").Value = '24.89'

foreach ($cell in $textForceCells) {
    $ws.Range($cell).Style = "Normal"
}

# --- Remaining cells: plain text updates (price strings with separators, percentages) ---
$ws.Range("D2").Value = '67.179.31'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '3.137.58'
$ws.Range("E3").Value = '  +3.43%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  +3.85%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '3.133.36'
$ws.Range("E8").Value = '  +3.38%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("E10").Value = '  -2.38%  '
$ws.Range("E11").Value = '  +1.96%  '
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("E13").Value = '  +0.67%  '
$ws.Range("E14").Value = '  +2.01%  '
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = '3.655.37'
$ws.Range("E16").Value = '  +3.39%  '
$ws.Range("D17").Value = '67.180.67'
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("D19").Value = '3.137.56'
$ws.Range("E19").Value = '  +3.52%  '
$ws.Range("E20").Value = '  -2.34%  '
$ws.Range("E21").Value = '  +4.56%  '
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("E23").Value = '  +3.72%  '
$ws.Range("E24").Value = '  +1.37%  '
$ws.Range("E25").Value = '  +4.15%  '
$ws.Range("E26").Value = '  +3.02%  '
$ws.Range("E27").Value = '  +0.24%  '
$ws.Range("E29").Value = '  -2.72%  '
$ws.Range("E31").Value = '  +1.21%  '
$ws.Range("E32").Value = '  +2.66%  '
$ws.Range("D33").Value = '0.0₃0997'
$ws.Range("E34").Value = '  -3.11%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("E36").Value = '  +1.45%  '
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("E39").Value = '  +2.76%  '
$ws.Range("E40").Value = '  +1.21%  '
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("E44").Value = '  -1.93%  '
$ws.Range("D45").Value = '2.845.88'
$ws.Range("E45").Value = '  +5.18%  '
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("E47").Value = '  +1.09%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("E50").Value = '  +1.71%  '
$ws.Range("E51").Value = '  -0.56%  '
